$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 41, shifting existing data (rows 41-58) down to (42-59).
$ws.Rows.Item(41).Insert()

# The new row 41 is populated as a duplicate of the row that is now directly below it
# (row 42, which holds the data that used to be row 41), then its date and volume are
# updated to the new weekly record's values.
$ws.Rows.Item(42).Copy()
$ws.Rows.Item(41).PasteSpecial()

$ws.Cells.Item(41, 4).Value = 44606
$ws.Cells.Item(41, 10).Value = 130
